$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A38").Value = 37
$ws.Range("C38").Value = "//a[contains(text(),'Privacy Policy')]/parent::span/parent::div"
$ws.Range("B38").Value = "BecomePartnerPage_Modal_Label_ViewPrivacyPolicy"

$ws.Range("B38").Select()
